# Adds a second line ("(TLVs)" / "TLVs") to the "Headers" textbox on slides
# 21, 22, 26, 27 and 28 (the textbox uses <a:spAutoFit/>, so PowerPoint grows
# its height automatically once the extra paragraph is added) and repositions
# (and, on slide 27, also resizes) the two "TLV description" textboxes on the
# right-hand side of slides 22, 26, 27 and 28 to match their new layout.
#
# Shape.Left/Top/Width/Height are COM `Single` (32-bit float) properties
# measured in points; PowerPoint stores shape geometry in EMU internally
# (1 pt = 12700 EMU). To land on an exact target EMU value despite the
# float32 round-trip, the point literals below were solved (in Python,
# using struct.pack('f', ...) to mimic the Single conversion) so that
# point_value * 12700, truncated, reproduces the exact target EMU.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 21 - "Headers" TextBox (shape 8) gains a second paragraph "(TLVs)"
# off stays (778877, 2053614) EMU; ext cy grows 307777 -> 523220 EMU via
# auto-fit, no other shapes on this slide change.
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$s21.Shapes.Item(8).TextFrame.TextRange.Text = "Headers`r(TLVs)"

# ---------------------------------------------------------------------------
# Slide 22 - "Headers" TextBox (shape 8) gains a second paragraph "TLVs";
# the two right-hand label textboxes (10, 12) are repositioned.
# ---------------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$s22.Shapes.Item(8).TextFrame.TextRange.Text = "Headers`rTLVs"

$sh = $s22.Shapes.Item(10)           # "Contiguous Label Stack ..." -> off 7275201,1064113 EMU
$sh.Left = 572.8505249023438
$sh.Top  = 83.7884292602539

$sh = $s22.Shapes.Item(12)           # "CW/G-ACH in TLV format ..." -> off 7283910,2030002 EMU
$sh.Left = 573.5362548828125
$sh.Top  = 159.84268188476562

# ---------------------------------------------------------------------------
# Slide 26 - "Headers" TextBox (shape 9) gains a second paragraph "TLVs";
# the two right-hand label textboxes (11, 13) are repositioned.
# ---------------------------------------------------------------------------
$s26 = $p.Slides.Item(26)
$s26.Shapes.Item(9).TextFrame.TextRange.Text = "Headers`rTLVs"

$sh = $s26.Shapes.Item(11)           # "Contiguous Label Stack ..." -> off 7456796,935223 EMU
$sh.Left = 587.1492919921875
$sh.Top  = 73.63961029052734

$sh = $s26.Shapes.Item(13)           # "CW/G-ACH in TLV format ..." -> off 7458497,1765882 EMU
$sh.Left = 587.2832641601562
$sh.Top  = 139.04583740234375

# ---------------------------------------------------------------------------
# Slide 27 - "Headers" TextBox (shape 8) gains a second paragraph "TLVs";
# the two right-hand label textboxes (11, 13) are repositioned, and the
# first one (11) is also resized narrower (cx 1603997 -> 1524000 EMU).
# ---------------------------------------------------------------------------
$s27 = $p.Slides.Item(27)
$s27.Shapes.Item(8).TextFrame.TextRange.Text = "Headers`rTLVs"

$sh = $s27.Shapes.Item(11)           # "Contiguous Label Stack ..." -> off 7581900,1113558 EMU, cx 1524000 EMU
$sh.Left  = 597.0
$sh.Top   = 87.6817398071289
$sh.Width = 120.0

$sh = $s27.Shapes.Item(13)           # "CW/G-ACH in TLV format ..." -> off 7581900,1975909 EMU
$sh.Left = 597.0
$sh.Top  = 155.58338928222656

# ---------------------------------------------------------------------------
# Slide 28 - "Headers" TextBox (shape 8) gains a second paragraph "TLVs";
# the two right-hand label textboxes (10, 12) are repositioned.
# ---------------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$s28.Shapes.Item(8).TextFrame.TextRange.Text = "Headers`rTLVs"

$sh = $s28.Shapes.Item(10)           # "Contiguous Label Stack ..." -> off 7084702,1161950 EMU
$sh.Left = 557.8505859375
$sh.Top  = 91.49212646484375

$sh = $s28.Shapes.Item(12)           # "CW/G-ACH in TLV format ..." -> off 7084702,2018335 EMU
$sh.Left = 557.8505859375
$sh.Top  = 158.9240264892578
